# Edit script for cmip6_miroc_responsible_parties.xlsx
# Adds 14 new responsible-party rows (14-27) to the "Responsibile Parties" sheet
# and updates the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Responsibile Parties")
$ws.Activate()

$ws.Range("A14").Value = "Kodama"
$ws.Range("B14").Value = "Chihiro Kodama"
$ws.Range("C14").Value = $false
$ws.Range("E14").Value = "kodamac@jamstec.go.jp"
$ws.Range("A15").Value = "Ohno"
$ws.Range("B15").Value = "Tomoki Ohno"
$ws.Range("C15").Value = $false
$ws.Range("E15").Value = "t-ohno@jamstec.go.jp"
$ws.Range("A16").Value = "Seiki"
$ws.Range("B16").Value = "Tatsuya Seiki"
$ws.Range("C16").Value = $false
$ws.Range("E16").Value = "tseiki@jamstec.go.jp"
$ws.Range("A17").Value = "Yashiro"
$ws.Range("B17").Value = "Hisashi Yashiro"
$ws.Range("C17").Value = $false
$ws.Range("E17").Value = "h.yashiro@riken.jp"
$ws.Range("A18").Value = "Noda"
$ws.Range("B18").Value = "Akira T. Noda"
$ws.Range("C18").Value = $false
$ws.Range("E18").Value = "a_noda@jamstec.go.jp"
$ws.Range("A19").Value = "Nakano"
$ws.Range("B19").Value = "Masuo Nakano"
$ws.Range("C19").Value = $false
$ws.Range("E19").Value = "masuo@jamstec.go.jp"
$ws.Range("A20").Value = "Yamada"
$ws.Range("B20").Value = "Yohei Yamada"
$ws.Range("C20").Value = $false
$ws.Range("E20").Value = "yoheiy@jamstec.go.jp"
$ws.Range("A21").Value = "Roh"
$ws.Range("B21").Value = "Woosub Roh"
$ws.Range("C21").Value = $false
$ws.Range("E21").Value = "ws-roh@aori.u-tokyo.ac.jp"
$ws.Range("A22").Value = "Satoh"
$ws.Range("B22").Value = "Masaki Satoh"
$ws.Range("C22").Value = $false
$ws.Range("E22").Value = "satoh@aori.u-tokyo.ac.jp"
$ws.Range("A23").Value = "Nitta"
$ws.Range("B23").Value = "Tomoko Nitta"
$ws.Range("C23").Value = $false
$ws.Range("E23").Value = "t_nitta@aori.u-tokyo.ac.jp"
$ws.Range("A24").Value = "Nasuno"
$ws.Range("B24").Value = "Tomoe Nasuno"
$ws.Range("C24").Value = $false
$ws.Range("E24").Value = "nasuno@jamstec.go.jp"
$ws.Range("A25").Value = "Miyakawa"
$ws.Range("B25").Value = "Tomoki Miyakawa"
$ws.Range("C25").Value = $false
$ws.Range("E25").Value = "miyakawa@aori.u-tokyo.ac.jp"
$ws.Range("A26").Value = "Chen"
$ws.Range("B26").Value = "Ying-Wen Chen"
$ws.Range("C26").Value = $false
$ws.Range("E26").Value = "yingwen@aori.u-tokyo.ac.jp"
$ws.Range("A27").Value = "Sugi"
$ws.Range("B27").Value = "Masato Sugi"
$ws.Range("C27").Value = $false
$ws.Range("E27").Value = "msugi@mri-jma.go.jp"

# Update the view's active cell / selection to match the edited area.
$ws.Range("C24").Select()
